$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update values only
$ws.Range("B3").Value = 2405465389304.486
$ws.Range("C3").Value = 2425524891177.107
$ws.Range("D3").Value = 665434837958216.5

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2500259559811.525
$ws.Range("C4").Value = 2435619873602.625
$ws.Range("D4").Value = 135167942608947.8

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 577425208653826.9
$ws.Range("C5").Value = 943062717515231.5
$ws.Range("D5").Value = 4391259312892450
